# Update "want to go" (想去人数) counts on several rows across sheets,
# matching the regenerated data snapshot at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 459
$ws1.Range("F10").Value = 159
$ws1.Range("F15").Value = 198
$ws1.Range("F16").Value = 1529
$ws1.Range("F19").Value = 357
$ws1.Range("F21").Value = 848
$ws1.Range("F25").Value = 2682
$ws1.Range("F30").Value = 643
$ws1.Range("F31").Value = 1322
$ws1.Range("F37").Value = 648

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 649

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F13").Value = 459
$ws4.Range("F17").Value = 159
$ws4.Range("F21").Value = 198
$ws4.Range("F22").Value = 1529
$ws4.Range("F25").Value = 357
$ws4.Range("F30").Value = 2682
$ws4.Range("F36").Value = 1322
$ws4.Range("F42").Value = 648
